$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 1
$ws.Range("C11").Value = 1
$ws.Range("C12").Value = 1
$ws.Range("C14").Value = 1
$ws.Range("C15").Value = 1
$ws.Range("C16").Value = 1
$ws.Range("C18").Value = 1
$ws.Range("C19").Value = 0
$ws.Range("C20").Value = 1
$ws.Range("C22").Value = 1
$ws.Range("C26").Value = 0
$ws.Range("C27").Value = 1
$ws.Range("C28").Value = 1
$ws.Range("C29").Value = 0
$ws.Range("C30").Value = 1
$ws.Range("C31").Value = 0
$ws.Range("C33").Value = 1
$ws.Range("C35").Value = 1
$ws.Range("C38").Value = 1
$ws.Range("C40").Value = 1
$ws.Range("C41").Value = 0
$ws.Range("C42").Value = 1
$ws.Range("C45").Value = 0
$ws.Range("C46").Value = 0
$ws.Range("C47").Value = 1
$ws.Range("C48").Value = 1
$ws.Range("C50").Value = 1

$excel.Calculate()
